$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section header (row 19) - shared string index 11 "3rd iteration - 60/20/20 split"
$ws.Range("A19").Value = "3rd iteration - 60/20/20 split"

# Column headers (row 20) reuse existing shared strings: Train, Val, Test, Total
$ws.Range("B20").Value = "Train"
$ws.Range("C20").Value = "Val"
$ws.Range("D20").Value = "Test"
$ws.Range("E20").Value = "Total"

# Row 21 - reef_1
$ws.Range("A21").Value = "reef_1"
$ws.Range("B21").Formula = "=296/2"
$ws.Range("C21").Formula = "=98/2"
$ws.Range("D21").Value = 51
$ws.Range("E21").Formula = "=SUM(B21:D21)"

# Row 22 - reef_2
$ws.Range("A22").Value = "reef_2"
$ws.Range("B22").Formula = "=308/2"
$ws.Range("C22").Formula = "=102/2"
$ws.Range("D22").Value = 52
$ws.Range("E22:E26").Formula = "=SUM(B22:D22)"

# Row 23 - reef_3
$ws.Range("A23").Value = "reef_3"
$ws.Range("B23").Formula = "=242/2"
$ws.Range("C23").Value = 40
$ws.Range("D23").Value = 41

# Row 24 - reef_4
$ws.Range("A24").Value = "reef_4"
$ws.Range("B24").Formula = "=244/2"
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 42

# Row 25 - reef_5
$ws.Range("A25").Value = "reef_5"
$ws.Range("B25").Formula = "=272/2"
$ws.Range("C25").Value = 45
$ws.Range("D25").Value = 47

# Row 26 - totals
$ws.Range("B26").Formula = "=SUM(B21:B25)"
$ws.Range("C26").Formula = "=SUM(C21:C25)"
$ws.Range("D26").Formula = "=SUM(D21:D25)"

# View state: scroll so row 3 is at top, select the final total cell
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C26").Select()
